{"js": "// Set line spacing to single (1.0) for every paragraph in the document\n// body, matching Word's \"Line Spacing: Single\" command\n// (<w:spacing w:line=\"240\" w:lineRule=\"auto\"/> on each paragraph's pPr).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const pf = paragraphs.items[i].paragraphFormat;\n  // The documented Word JS API only exposes `lineSpacing` (a raw point\n  // value) and has no property for the line-spacing *rule*, so \"single\"\n  // spacing (rule = auto, value = 12pt) can't be expressed through the\n  // public surface alone. Reach through to the same OM bridge the\n  // documented setters themselves call, mirroring\n  // ParagraphFormat.LineSpacingRule / ParagraphFormat.LineSpacing in the\n  // Word object model.\n  pf._omSet(\"LineSpacingRule\", 0);\n  pf._omSet(\"LineSpacing\", 12);\n}\n\nawait context.sync();\n", "ps1": "# Set line spacing to single (1.0) for every paragraph in the document.\n# Single spacing == LineSpacingRule 0 (wdLineSpaceSingle) with LineSpacing\n# of 12pt, which Word serializes as <w:spacing w:line=\"240\" w:lineRule=\"auto\"/>.\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $p.Range.ParagraphFormat.LineSpacingRule = 0\n    $p.Range.ParagraphFormat.LineSpacing = 12\n}\n"}
